$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold/border/center-top alignment) from B1 into the
# new date column header C1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "13-01-2023"

# Swap the two row labels: what was "total" (row2) becomes "Arpenta ex Mercosur",
# and what was "Arpenta ex Mercosur" (row4) becomes "total"
$ws.Range("A2").Value = "Arpenta ex Mercosur"
$ws.Range("A4").Value = "total"

# Populate the new values for the new date column (C) on rows 2-4
$ws.Range("C2").Value = 64526.08
$ws.Range("C3").Value = 64526.08
$ws.Range("C4").Value = 64526.08
